$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only the specific "Price" column cells whose values change need to be force-typed as
# text (NumberFormat "@") before the assignment below, so that numeric-looking strings
# such as "0.688" or "56.60" are preserved verbatim instead of being parsed into floats
# by Excel. Cells that are not touched by this update keep their original formatting.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '35.846.83'
$ws.Range('E2').Value = '  +1.39%  '
$ws.Range('D3').Value = '1.899.16'
$ws.Range('E3').Value = '  +0.78%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '247.22'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').Value = '0.688'
$ws.Range('E6').Value = '  +0.34%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '42.93'
$ws.Range('E8').Value = '  +0.86%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = '0.363'
$ws.Range('E9').Value = '  +3.14%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').Value = '56.60'
$ws.Range('E10').Value = '  +7.28%  '
$ws.Range('E11').Value = '  +3.16%  '
$ws.Range('D12').Value = '0.0988'
$ws.Range('E12').Value = '  +2.06%  '
$ws.Range('D13').Value = '14.83'
$ws.Range('E13').Value = '  +14.45%  '
$ws.Range('D14').Value = '0.793'
$ws.Range('E14').Value = '  +7.55%  '
$ws.Range('D15').Value = '2.172.97'
$ws.Range('E15').Value = '  +0.75%  '
$ws.Range('E16').Value = '  +3.51%  '
$ws.Range('D17').Value = '1.899.46'
$ws.Range('E17').Value = '  +0.93%  '
$ws.Range('D18').Value = '35.868.44'
$ws.Range('E18').Value = '  +1.38%  '
$ws.Range('D19').Value = '73.40'
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').Value = '0.0₃0832'
$ws.Range('E20').Value = '  +1.49%  '
$ws.Range('D21').Value = '247.19'
$ws.Range('E21').Value = '  +1.44%  '
$ws.Range('D22').Value = '13.04'
$ws.Range('E22').Value = '  +2.43%  '
$ws.Range('D23').Value = '5.21'
$ws.Range('E23').Value = '  +5.91%  '
$ws.Range('D24').Value = '2.67'
$ws.Range('E24').Value = '  +3.70%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').Value = '2.21'
$ws.Range('E26').Value = '  +0.96%  '
$ws.Range('D27').Value = '166.96'
$ws.Range('E27').Value = '  +0.99%  '
$ws.Range('D28').Value = '8.75'
$ws.Range('E28').Value = '  +3.89%  '
$ws.Range('D29').Value = '18.44'
$ws.Range('E29').Value = '  +0.56%  '
$ws.Range('E30').Value = '  +1.30%  '
$ws.Range('D31').Value = '4.50'
$ws.Range('E31').Value = '  +6.28%  '
$ws.Range('D32').Value = '0.0609'
$ws.Range('E32').Value = '  +5.94%  '
$ws.Range('D33').Value = '4.32'
$ws.Range('E33').Value = '  +4.35%  '
$ws.Range('D34').Value = '1.92'
$ws.Range('E34').Value = '  +3.01%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Value = '1.47'
$ws.Range('E36').Value = '  -14.71%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = '0.0799'
$ws.Range('E37').Value = '  +17.26%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '0.852'
$ws.Range('E38').Value = '  +1.74%  '
$ws.Range('D39').Value = '1.97'
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('D40').Value = '0.0227'
$ws.Range('E40').Value = '  +2.88%  '
$ws.Range('D41').Value = '99.31'
$ws.Range('E41').Value = '  +3.01%  '
$ws.Range('D42').Value = '15.04'
$ws.Range('E42').Value = '  +22.15%  '
$ws.Range('D43').Value = '16.65'
$ws.Range('E43').Value = '  -1.50%  '
$ws.Range('D44').Value = '1.09'
$ws.Range('E44').Value = '  +1.56%  '
$ws.Range('D45').Value = '1.321.64'
$ws.Range('E45').Value = '  +2.53%  '
$ws.Range('E46').Value = '  +1.40%  '
$ws.Range('D47').Value = '0.0808'
$ws.Range('E47').Value = '  +1.58%  '
$ws.Range('E48').Value = '  +0.56%  '
$ws.Range('D49').Value = '2.77'
$ws.Range('E49').Value = '  +1.16%  '
$ws.Range('D50').Value = '6.34'
$ws.Range('E50').Value = '  +1.70%  '
$ws.Range('D51').Value = '42.51'
$ws.Range('E51').Value = '  +0.10%  '
